# Weekly roll-forward of the "Agricola del Norte S.A. de Arica - Pina" sheet.
#
# The whole historical block of rows (125-215) shifts down by exactly one
# week-block (4 rows): new row R (129<=R<=219) takes on the content that
# used to live at row R-4. That is accomplished with a single range copy
# (Excel recalculates dates/styles/values together, exactly like a manual
# "select block, copy, paste 4 rows down").
#
# The vacated block at the top (rows 125-128) is then populated with the
# brand-new week's figures (fecha 2022-08-17, serial 44790).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole trailing block (rows 125-215) down by 4 rows, in one go,
# which both fills rows 129-219 (including the brand new rows 216-219) and
# carries formatting (e.g. the date number format on column D) along.
$ws.Range("A125:T215").Copy($ws.Range("A129:T219"))

# Now populate the newly-freed rows 125-128 with this week's data.

# Row 125: Especial
$ws.Range("D125").Value = 44790
$ws.Range("N125").Value = 20000
$ws.Range("O125").Value = 21000
$ws.Range("P125").Value = 20500
$ws.Range("S125").Value = 2050

# Row 126: Primera
$ws.Range("D126").Value = 44790
$ws.Range("M126").Value = 250
$ws.Range("N126").Value = 20000
$ws.Range("O126").Value = 21000
$ws.Range("P126").Value = 20500
$ws.Range("S126").Value = 1708

# Row 127: Segunda
$ws.Range("D127").Value = 44790
$ws.Range("N127").Value = 20000
$ws.Range("O127").Value = 21000
$ws.Range("P127").Value = 20500
$ws.Range("S127").Value = 1464

# Row 128: Tercera
$ws.Range("D128").Value = 44790
$ws.Range("M128").Value = 250
$ws.Range("N128").Value = 20000
$ws.Range("O128").Value = 21000
$ws.Range("P128").Value = 20500
$ws.Range("S128").Value = 1281
